$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swap: 25 <-> 26 (Toncoin / PancakeSwap reorder) ---
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.35'
$ws.Range("E26").Value = '  +16.24%  '

# --- Row swap: 40 <-> 41 (Kaspa / ARBITRUM reorder) ---
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0674'
$ws.Range("E40").Value = '  +14.18%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.11'
$ws.Range("E41").Value = '  +0.52%  '

# --- Price / volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.269.90'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.25'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.69'
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.692'
$ws.Range("E6").Value = '  +9.64%  '
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.82'
$ws.Range("E8").Value = '  -3.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.351'
$ws.Range("E9").Value = '  +5.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '54.07'
$ws.Range("E10").Value = '  +14.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0729'
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.181.42'
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.38'
$ws.Range("E14").Value = '  +4.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.706'
$ws.Range("E15").Value = '  +2.97%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.907.94'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.279.60'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.46'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0824'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.49'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("E22").Value = '  +1.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.25'
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.59'
$ws.Range("E28").Value = '  +4.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.44'
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.130'
$ws.Range("E30").Value = '  +3.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.977'
$ws.Range("E32").Value = '  +9.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.19'
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  +0.94%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0211'
$ws.Range("E42").Value = '  +3.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.25'
$ws.Range("E43").Value = '  +6.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '90.79'
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.346.20'
$ws.Range("E46").Value = '  +2.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.02'
$ws.Range("E47").Value = '  +2.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.53'
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.79'
$ws.Range("E49").Value = '  +2.25%  '
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.56'
$ws.Range("E51").Value = '  -2.60%  '
